$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'99.264.16"
$ws.Range("E2").Value = "  +0.89%  "

# Row 3
$ws.Range("D3").Value = "'3.290.59"
$ws.Range("E3").Value = "  -2.51%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").Value = "'254.37"
$ws.Range("E5").Value = "  -1.35%  "

# Row 6
$ws.Range("D6").Value = "'622.10"
$ws.Range("E6").Value = "  -0.28%  "

# Row 7
$ws.Range("D7").Value = "'1.42"
$ws.Range("E7").Value = "  +15.69%  "

# Row 8
$ws.Range("D8").Value = "'0.401"
$ws.Range("E8").Value = "  +3.77%  "

# Row 9
$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "  -0.01%  "

# Row 10
$ws.Range("D10").Value = "'0.958"
$ws.Range("E10").Value = "  +19.25%  "

# Row 11
$ws.Range("D11").Value = "'3.286.52"
$ws.Range("E11").Value = "  -2.50%  "

# Row 12
$ws.Range("D12").Value = "'0.199"
$ws.Range("E12").Value = "  -0.36%  "

# Row 13
$ws.Range("D13").Value = "'39.26"
$ws.Range("E13").Value = "  +8.93%  "

# Row 14
$ws.Range("D14").Value = "'98.925.28"
$ws.Range("E14").Value = "  +0.99%  "

# Row 15
$ws.Range("E15").Value = "  +0.07%  "

# Row 16
$ws.Range("D16").Value = "'3.882.20"
$ws.Range("E16").Value = "  -2.95%  "

# Row 17
$ws.Range("D17").Value = "'5.45"
$ws.Range("E17").Value = "  -1.06%  "

# Row 18
$ws.Range("D18").Value = "'3.275.51"
$ws.Range("E18").Value = "  -2.82%  "

# Row 19
$ws.Range("D19").Value = "'3.44"
$ws.Range("E19").Value = "  -4.37%  "

# Row 20
$ws.Range("D20").Value = "'15.31"
$ws.Range("E20").Value = "  +1.33%  "

# Row 21
$ws.Range("D21").Value = "'6.30"
$ws.Range("E21").Value = "  +7.36%  "

# Row 22
$ws.Range("D22").Value = "'484.88"
$ws.Range("E22").Value = "  +0.01%  "

# Row 23
$ws.Range("E23").Value = "  +0.29%  "

# Row 24
$ws.Range("E24").Value = "  -3.86%  "

# Row 25
$ws.Range("D25").Value = "'5.61"
$ws.Range("E25").Value = "  -1.30%  "

# Row 26
$ws.Range("D26").Value = "'88.54"
$ws.Range("E26").Value = "  +0.12%  "

# Row 27
$ws.Range("D27").Value = "'0.316"
$ws.Range("E27").Value = "  +22.55%  "

# Row 28
$ws.Range("D28").Value = "'11.96"
$ws.Range("E28").Value = "  -0.66%  "

# Row 29
$ws.Range("D29").Value = "'3.431.03"
$ws.Range("E29").Value = "  -3.52%  "

# Row 30
$ws.Range("E30").Value = "  -0.06%  "

# Row 31
$ws.Range("B31").Value = "Cronos"
$ws.Range("C31").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D31").Value = "'0.189"
$ws.Range("E31").Value = "  +2.20%  "

# Row 32
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "'0.135"
$ws.Range("E32").Value = "  +7.77%  "

# Row 33
$ws.Range("D33").Value = "'10.29"
$ws.Range("E33").Value = "  +10.18%  "

# Row 34
$ws.Range("D34").Value = "'0.998"
$ws.Range("E34").Value = "  -0.09%  "

# Row 35
$ws.Range("D35").Value = "'27.80"
$ws.Range("E35").Value = "  +1.68%  "

# Row 36
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "'0.148"
$ws.Range("E36").Value = "  -1.65%  "

# Row 37
$ws.Range("B37").Value = "PolygonEcosystemToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D37").Value = "'0.471"
$ws.Range("E37").Value = "  +5.32%  "

# Row 38
$ws.Range("D38").Value = "'7.19"
$ws.Range("E38").Value = "  -3.46%  "

# Row 39
$ws.Range("E39").Value = "  -0.77%  "

# Row 40
$ws.Range("D40").Value = "'24.84"
$ws.Range("E40").Value = "  -0.26%  "

# Row 41
$ws.Range("D41").Value = "'489.34"
$ws.Range("E41").Value = "  -5.86%  "

# Row 42
$ws.Range("B42").Value = "Fetch.AI"
$ws.Range("C42").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D42").Value = "'1.23"
$ws.Range("E42").Value = "  -2.77%  "

# Row 43
$ws.Range("B43").Value = "MantraDAO"
$ws.Range("C43").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D43").Value = "'3.60"
$ws.Range("E43").Value = "  -0.42%  "

# Row 44
$ws.Range("E44").Value = "  -0.02%  "

# Row 45
$ws.Range("D45").Value = "'0.769"
$ws.Range("E45").Value = "  -0.96%  "

# Row 46
$ws.Range("E46").Value = "  -6.00%  "

# Row 47
$ws.Range("D47").Value = "'1.95"
$ws.Range("E47").Value = "  +1.36%  "

# Row 48
$ws.Range("D48").Value = "'157.56"
$ws.Range("E48").Value = "  -2.09%  "

# Row 49
$ws.Range("D49").Value = "'0.843"
$ws.Range("E49").Value = "  +6.08%  "

# Row 50
$ws.Range("D50").Value = "'7.23"
$ws.Range("E50").Value = "  +14.21%  "

# Row 51
$ws.Range("E51").Value = "  +3.60%  "

Write-Output "done"